$wb = $excel.ActiveWorkbook

# 1) Duplicate the "MF Requirements" sheet (same 22-column schema/style) and place the
#    copy immediately before "CbCR Notifications" -> this becomes the new "LF Requirements" tab.
$sourceSheet = $wb.Worksheets.Item("MF Requirements")
$beforeSheet = $wb.Worksheets.Item("CbCR Notifications")
$sourceSheet.Copy($beforeSheet)
$ws = $wb.Worksheets.Item("MF Requirements (2)")
$ws.Name = "LF Requirements"

# 2) The source sheet only has 5 data rows (2-6); LF Requirements needs 6 (2-7), so
#    clone the formatting (fill/font) of the last data row down into row 7 first.
$ws.Range("A6:V6").Copy()
$ws.Range("A7:V7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# "Condition Group" (col H) holds short codes like "1"/"2" and must stay text, not numeric.
$ws.Range("H2:H7").NumberFormat = "@"

# Row 2: Germany
$ws.Range("A2").Value = "Germany"
$ws.Range("B2").Value = "Conditional"
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = "No"
$ws.Range("G2").Value = "LF-DE-1"
$ws.Range("H2").Value = "1"
$ws.Range("I2").Value = "OR"
$ws.Range("J2").Value = "RPTs"
$ws.Range("K2").Value = "Transaction (Goods)"
$ws.Range("L2").Value = 6000000
$ws.Range("M2").Value = "EUR"
$ws.Range("N2").Value = ">"
$ws.Range("O2").Value = "CIT Date"
$ws.Range("P2").Value = "Expected 31 Jul (CIT filing)"
$ws.Range("Q2").Value = "Upon Request"
$ws.Range("R2").Value = "Within 30 days of audit notice"
$ws.Range("S2").Value = 30
$ws.Range("T2").Value = "FY2024"
$ws.Range("U2").Value = "LF required if goods RPTs exceed 6M EUR"
$ws.Range("V2").Value = "Automatic submission upon audit"

# Row 3: Germany
$ws.Range("A3").Value = "Germany"
$ws.Range("B3").Value = "Conditional"
$ws.Range("C3").Value = ""
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = ""
$ws.Range("F3").Value = "No"
$ws.Range("G3").Value = "LF-DE-1"
$ws.Range("H3").Value = "2"
$ws.Range("I3").Value = "OR"
$ws.Range("J3").Value = "RPTs"
$ws.Range("K3").Value = "Transaction (Services)"
$ws.Range("L3").Value = 600000
$ws.Range("M3").Value = "EUR"
$ws.Range("N3").Value = ">"
$ws.Range("O3").Value = "CIT Date"
$ws.Range("P3").Value = "Expected 31 Jul (CIT filing)"
$ws.Range("Q3").Value = "Upon Request"
$ws.Range("R3").Value = "Within 30 days of audit notice"
$ws.Range("S3").Value = 30
$ws.Range("T3").Value = "FY2024"
$ws.Range("U3").Value = "OR services/other RPTs exceed 600K EUR"
$ws.Range("V3").Value = "Automatic submission upon audit"

# Row 4: Spain
$ws.Range("A4").Value = "Spain"
$ws.Range("B4").Value = "Conditional"
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = ""
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = "No"
$ws.Range("G4").Value = "LF-ES-1"
$ws.Range("H4").Value = "1"
$ws.Range("I4").Value = "OR"
$ws.Range("J4").Value = "RPTs"
$ws.Range("K4").Value = "Transaction (All)"
$ws.Range("L4").Value = 250000
$ws.Range("M4").Value = "EUR"
$ws.Range("N4").Value = ">"
$ws.Range("O4").Value = "CIT Date"
$ws.Range("P4").Value = "Expected 25 Jul"
$ws.Range("Q4").Value = "Upon Request"
$ws.Range("R4").Value = "Within 10 days of request"
$ws.Range("S4").Value = 10
$ws.Range("T4").Value = "FY2016"
$ws.Range("U4").Value = "LF required if local RPTs exceed 250K EUR"
$ws.Range("V4").Value = "Maintain contemporaneously"

# Row 5: Malaysia
$ws.Range("A5").Value = "Malaysia"
$ws.Range("B5").Value = "Always"
$ws.Range("C5").Value = ""
$ws.Range("D5").Value = ""
$ws.Range("E5").Value = ""
$ws.Range("F5").Value = "No"
$ws.Range("G5").Value = ""
$ws.Range("H5").Value = ""
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = ""
$ws.Range("L5").Value = ""
$ws.Range("M5").Value = ""
$ws.Range("N5").Value = ""
$ws.Range("O5").Value = "CIT Date"
$ws.Range("P5").Value = "By 7 months after FYE (CIT filing)"
$ws.Range("Q5").Value = "Upon Request"
$ws.Range("R5").Value = "Within 14 days"
$ws.Range("S5").Value = 14
$ws.Range("T5").Value = "FY2023"
$ws.Range("U5").Value = "CTPD (LF) required for all entities with RPTs. MF content integrated per 2023 TPD."
$ws.Range("V5").Value = "File with CIT return"

# Row 6: United States
$ws.Range("A6").Value = "United States"
$ws.Range("B6").Value = "Conditional"
$ws.Range("C6").Value = ""
$ws.Range("D6").Value = ""
$ws.Range("E6").Value = ""
$ws.Range("F6").Value = "Yes"
$ws.Range("G6").Value = ""
$ws.Range("H6").Value = ""
$ws.Range("I6").Value = ""
$ws.Range("J6").Value = ""
$ws.Range("K6").Value = ""
$ws.Range("L6").Value = ""
$ws.Range("M6").Value = ""
$ws.Range("N6").Value = ""
$ws.Range("O6").Value = "None"
$ws.Range("P6").Value = "Voluntary preparation recommended"
$ws.Range("Q6").Value = "None"
$ws.Range("R6").Value = "N/A - voluntary"
$ws.Range("S6").Value = ""
$ws.Range("T6").Value = "FY2018"
$ws.Range("U6").Value = "Voluntary LF preparation for penalty protection under IRC §6662. No filing requirement."
$ws.Range("V6").Value = "Contemporaneous documentation provides reasonable cause defense"

# Row 7: Canada
$ws.Range("A7").Value = "Canada"
$ws.Range("B7").Value = "Conditional"
$ws.Range("C7").Value = ""
$ws.Range("D7").Value = ""
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = "Yes"
$ws.Range("G7").Value = ""
$ws.Range("H7").Value = ""
$ws.Range("I7").Value = ""
$ws.Range("J7").Value = ""
$ws.Range("K7").Value = ""
$ws.Range("L7").Value = ""
$ws.Range("M7").Value = ""
$ws.Range("N7").Value = ""
$ws.Range("O7").Value = "None"
$ws.Range("P7").Value = "Voluntary preparation recommended"
$ws.Range("Q7").Value = "None"
$ws.Range("R7").Value = "N/A - voluntary"
$ws.Range("S7").Value = ""
$ws.Range("T7").Value = "FY2015"
$ws.Range("U7").Value = "Voluntary LF for penalty protection. No statutory filing requirement."
$ws.Range("V7").Value = "Contemporaneous documentation required for transfer pricing adjustment defense"

# 3) Restore the originally-active sheet/selection so we don't leave an unrelated
#    "last viewed tab" change behind.
$excel.CutCopyMode = $false
$wb.Worksheets.Item("Data Dictionary").Activate()
[void]$wb.Worksheets.Item("Data Dictionary").Range("A1").Select()
